$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the "Existing Liability w/Notice Number" PaymentType label to the
#    new "Existing Liability with Notice/Invoice Number" wording, for all the
#    rows that used the old label (rows 2-13 and 34-41 in column D).
$ws.Range("D2:D13").Value = "Existing Liability with Notice/Invoice Number"
$ws.Range("D34:D41").Value = "Existing Liability with Notice/Invoice Number"

# 2) Mark the "Execute" column (C) as "Y" for all the rows that previously had
#    no value there (rows 2-20 and 28-54). Rows 21-27 already had "Y" and are
#    left untouched.
$ws.Range("C2:C20").Value = "Y"
$ws.Range("C28:C54").Value = "Y"

# 3) The wider new PaymentType text no longer fits the old column D width, so
#    re-fit column D to the new content (mirrors what Excel does automatically
#    for the author after retyping the longer label).
$ws.Columns.Item(4).AutoFit()

# 4) Reflect the author's final selection/cursor position when the workbook
#    was saved.
$ws.Range("C24:C54").Select()
